$d = $word.ActiveDocument

function Find-Range([string]$text) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Find failed for: $text"
    }
    return $d.Range($rng.Start, $rng.End)
}

function Replace-Plain([string]$find, [string]$replaceWith) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replaceWith, 2)
    if (-not $ok) {
        throw "Replace failed for: $find"
    }
}

function Insert-RawXml([string]$find, [string]$bodyXml) {
    $rng = Find-Range $find
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    [void]$rng.InsertXML($pkg)
}

# 1. "Capstone Project I" -> "Capstone Project II" (new trailing run) plus a
#    fresh "_GoBack" bookmark right after it, inside the same paragraph.
Insert-RawXml "Capstone Project I" (
    '<w:p w14:paraId="193CFFE4" w14:textId="77777777" w:rsidR="00782380" w:rsidRPr="00A245F2" w:rsidRDefault="00782380">' +
    '<w:pPr><w:pStyle w:val="Standard1"/><w:spacing w:before="0" w:after="0"/><w:rPr><w:b/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr>' +
    '<w:bookmarkStart w:id="2" w:name="Logistics"/><w:bookmarkEnd w:id="2"/>' +
    '<w:r w:rsidRPr="00A245F2"><w:rPr><w:b/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>Capstone Project I</w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>I</w:t></w:r>' +
    '<w:bookmarkStart w:id="3" w:name="_GoBack"/><w:bookmarkEnd w:id="3"/>' +
    '</w:p>'
)

# 2. "Arsalan Farooqui" + ": " + "Implement ..." collapse into a single run.
Replace-Plain "Arsalan Farooqui: Implement basic navigation structure in React Native for application" `
              "Arsalan Farooqui: Implement basic navigation structure in React Native for application"

# 3. "Guiseppe Ragusa" + ": " + "Begin learning ..." -> "G" + "iu" + "seppe Ragusa: Begin learning ..."
Insert-RawXml "Guiseppe Ragusa: Begin learning React library and make starter React project" (
    '<w:p w14:paraId="5D799F9D" w14:textId="77777777" w:rsidR="002D55F7" w:rsidRDefault="002D55F7" w:rsidP="002D55F7">' +
    '<w:pPr><w:pStyle w:val="Standard1"/><w:tabs><w:tab w:val="left" w:pos="780"/></w:tabs><w:spacing w:before="120" w:after="120"/></w:pPr>' +
    '<w:r><w:t>G</w:t></w:r><w:r><w:t>iu</w:t></w:r><w:r><w:t>seppe Ragusa: Begin learning React library and make starter React project</w:t></w:r>' +
    '</w:p>'
)

# 4. "Andrew Cobb : " + "Go through ..." collapse into a single run.
Replace-Plain "Andrew Cobb : Go through React Native CRUD tutorial" `
              "Andrew Cobb : Go through React Native CRUD tutorial"

# 5. "Implement basic react application" -> 5 runs describing the week 9 task.
Insert-RawXml "Implement basic react application" (
    '<w:p w14:paraId="7863403D" w14:textId="486D9DAF" w:rsidR="00356405" w:rsidRDefault="00356405" w:rsidP="00356405">' +
    '<w:pPr><w:pStyle w:val="Standard1"/><w:spacing w:before="120" w:after="120"/></w:pPr>' +
    '<w:r><w:t>Install</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> react </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">native </w:t></w:r>' +
    '<w:r><w:t>application</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> on local machine, begin implementing screen layouts</w:t></w:r>' +
    '</w:p>'
)

# 6. Second "Guiseppe Ragusa" (assignee cell, week 9 row) -> "G" + "iu" + "seppe Ragusa"
Insert-RawXml "Guiseppe Ragusa" (
    '<w:p w14:paraId="12A46D48" w14:textId="2CE1917B" w:rsidR="00356405" w:rsidRDefault="00356405" w:rsidP="00356405">' +
    '<w:pPr><w:pStyle w:val="Standard1"/><w:spacing w:before="120" w:after="120"/></w:pPr>' +
    '<w:r><w:t>G</w:t></w:r><w:r><w:t>iu</w:t></w:r><w:r><w:t>seppe Ragusa</w:t></w:r>' +
    '</w:p>'
)

# 7. Drop the stray "_GoBack" bookmark that used to sit after "create boat layouts page".
Insert-RawXml "create boat layouts page" (
    '<w:p w14:paraId="6241F3B4" w14:textId="6152CD40" w:rsidR="00356405" w:rsidRDefault="00356405" w:rsidP="00356405">' +
    '<w:pPr><w:pStyle w:val="Standard1"/><w:spacing w:before="120" w:after="120"/></w:pPr>' +
    '<w:r><w:t>Application routing and</w:t></w:r>' +
    '<w:r w:rsidR="00803CCF"><w:t xml:space="preserve"> create boat layouts page</w:t></w:r>' +
    '</w:p>'
)

# 8. "Date: Monday, January 20, 2020" -> 3 runs splitting out the "1" in "21".
Insert-RawXml "Date: Monday, January 20, 2020" (
    '<w:p w14:paraId="52F52F9A" w14:textId="39B15F14" w:rsidR="00356405" w:rsidRDefault="00356405" w:rsidP="00356405">' +
    '<w:pPr><w:pStyle w:val="Standard1"/><w:spacing w:before="120" w:after="120"/></w:pPr>' +
    '<w:r><w:t>Date: Monday, January 2</w:t></w:r>' +
    '<w:r><w:t>1</w:t></w:r>' +
    '<w:r><w:t>, 2020</w:t></w:r>' +
    '</w:p>'
)

# 9. Room "E430" -> "C410" (keep it as its own run, do not merge with "Campus Room ")
Insert-RawXml "Location: George Brown College Casa Loma   Campus Room E430" (
    '<w:p w14:paraId="06808F54" w14:textId="5201E30C" w:rsidR="00356405" w:rsidRDefault="00356405" w:rsidP="00356405">' +
    '<w:pPr><w:pStyle w:val="Standard1"/><w:spacing w:before="120" w:after="120"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Location: </w:t></w:r>' +
    '<w:r w:rsidRPr="002E154B"><w:t xml:space="preserve">George Brown College Casa Loma </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">  </w:t></w:r>' +
    '<w:r w:rsidRPr="002E154B"><w:t xml:space="preserve">Campus Room </w:t></w:r>' +
    '<w:r><w:t>C410</w:t></w:r>' +
    '</w:p>'
)

# 10. Append ".ca" (matching run formatting) after the truncated e-mail address.
Insert-RawXml "Team member 3 Giuseppe.Ragusa@georgebrown" (
    '<w:p w14:paraId="4C905EC5" w14:textId="08D35F3C" w:rsidR="00356405" w:rsidRPr="0016291E" w:rsidRDefault="00356405" w:rsidP="00356405">' +
    '<w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr>' +
    '<w:r w:rsidRPr="009D6683"><w:rPr><w:sz w:val="19"/><w:szCs w:val="19"/><w:lang w:val="en-GB"/></w:rPr><w:t>Team member 3 Giuseppe.Ragusa@georgebrown</w:t></w:r>' +
    '<w:r><w:rPr><w:sz w:val="19"/><w:szCs w:val="19"/><w:lang w:val="en-GB"/></w:rPr><w:t>.ca</w:t></w:r>' +
    '</w:p>'
)

Write-Output "done"
